# Update the "取得日時" (retrieved datetime) column (A) for the data rows
# on the "ランサーズ" sheet from 2026-02-02 07:00:27 to 2026-02-02 07:06:54.
# These cells hold plain text timestamps (no date number format applied),
# so assigning a string keeps them as text rather than a date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-02 07:06:54"

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
